$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old 2-column admin/password rows (A1:B6)
$ws.Range("A1:B6").Clear()

# New 3-column test-case table: testcase name, username, password
$ws.Range("A1").Value = "stestcasename"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "password"

$ws.Range("A2").Value = "AmazonLoginPageTest"
$ws.Range("B2").Value = "admin@abc.com"
$ws.Range("C2").Value = "admin"

$ws.Range("A3").Value = "OrangeHRMLoginTest"
$ws.Range("B3").Value = "Admin"
$ws.Range("C3").Value = "admin@123"

# Re-select and auto-fit the new layout
$ws.Range("A3").Select()
$ws.Columns("A:C").AutoFit()
